$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "eval" column (H) text for the rows whose assignments now
# include a Peer Review step.

$ws.Range("H3").Value = "Blog about data management preparation (Due Thu: 8/31) `n* Peer Review of DM blog post (Due 9/4)`n* Data management code file (Due 9/4)`n* Citation [Assignment](hw/Citation_Assignment.html) (Due 9/7)"

$ws.Range("H4").Value = "Univariate graphing assignment (Due 9/13)`n* Peer Revew (Due 9/18)`n"

$ws.Range("H5").Value = "Research plan outline (Due 9/18)`n* Peer Review (Due 9/20)"

$ws.Range("H7").Value = "Bivariate graphing assignment (Due 9/25)`n* Peer Review (Due 10/2)"

# Move the view: scroll so row 3 is at the top, and select E3.
$ws.Range("E3").Select()
$excel.ActiveWindow.ScrollRow = 3
